$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $r = $ws.Range($cell)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = $origStyle
}

# Row 2
$ws.Range("D2").Value = '25.647.08'
$ws.Range("E2").Value = '  -4.22%  '

# Row 3
$ws.Range("D3").Value = '1.813.02'
$ws.Range("E3").Value = '  -3.06%  '

# Row 4
Set-TextValue "D4" '0.9982'
$ws.Range("E4").Value = '  -0.25%  '

# Row 5
Set-TextValue "D5" '274.94'
$ws.Range("E5").Value = '  -8.39%  '

# Row 6
Set-TextValue "D6" '0.9983'
$ws.Range("E6").Value = '  -0.23%  '

# Row 7
Set-TextValue "D7" '0.5038'
$ws.Range("E7").Value = '  -5.40%  '

# Row 8
Set-TextValue "D8" '0.3429'
$ws.Range("E8").Value = '  -8.05%  '

# Row 9
Set-TextValue "D9" '44.04'
$ws.Range("E9").Value = '  -2.98%  '

# Row 10
Set-TextValue "D10" '0.06664'
$ws.Range("E10").Value = '  -6.87%  '

# Row 11
Set-TextValue "D11" '19.57'
$ws.Range("E11").Value = '  -8.84%  '

# Row 12
Set-TextValue "D12" '0.8005'
$ws.Range("E12").Value = '  -9.67%  '

# Row 13
Set-TextValue "D13" '0.07840'
$ws.Range("E13").Value = '  -3.88%  '

# Row 14
$ws.Range("D14").Value = '1.816.69'

# Row 15
Set-TextValue "D15" '5.036'
$ws.Range("E15").Value = '  -4.84%  '

# Row 16
Set-TextValue "D16" '87.37'
$ws.Range("E16").Value = '  -5.52%  '

# Row 17
Set-TextValue "D17" '0.9971'
$ws.Range("E17").Value = '  -0.32%  '

# Row 18
Set-TextValue "D18" '14.00'
$ws.Range("E18").Value = '  -5.65%  '

# Row 19
$ws.Range("E19").Value = '  +0.05%  '

# Row 20
Set-TextValue "D20" '0.000007978'
$ws.Range("E20").Value = '  -6.04%  '

# Row 21
$ws.Range("D21").Value = '25.660.75'
$ws.Range("E21").Value = '  -4.26%  '

# Row 22
Set-TextValue "D22" '4.718'
$ws.Range("E22").Value = '  -5.10%  '

# Row 23
Set-TextValue "D23" '9.908'
$ws.Range("E23").Value = '  -6.73%  '

# Row 24
Set-TextValue "D24" '6.118'
$ws.Range("E24").Value = '  -4.06%  '

# Row 25
Set-TextValue "D25" '2.262'
$ws.Range("E25").Value = '  -1.03%  '

# Row 26
Set-TextValue "D26" '142.42'
$ws.Range("E26").Value = '  -2.35%  '

# Row 27
Set-TextValue "D27" '1.655'
$ws.Range("E27").Value = '  -4.35%  '

# Row 28
Set-TextValue "D28" '17.07'
$ws.Range("E28").Value = '  -5.29%  '

# Row 29
Set-TextValue "D29" '108.56'
$ws.Range("E29").Value = '  -4.52%  '

# Row 30
Set-TextValue "D30" '4.269'
$ws.Range("E30").Value = '  -9.13%  '

# Row 31
Set-TextValue "D31" '4.211'
$ws.Range("E31").Value = '  -9.07%  '

# Row 32
Set-TextValue "D32" '0.08714'
$ws.Range("E32").Value = '  -4.29%  '

# Row 33
Set-TextValue "D33" '0.04782'
$ws.Range("E33").Value = '  -4.71%  '

# Row 34
Set-TextValue "D34" '1.132'
$ws.Range("E34").Value = '  -3.55%  '

# Row 35: 'HuobiToken' -> 'ImmutableX'
$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue "D35" '0.7189'
$ws.Range("E35").Value = '  -11.37%  '

# Row 36: 'ImmutableX' -> 'HuobiToken'
$ws.Range("B36").Value = 'HuobiToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue "D36" '2.835'
$ws.Range("E36").Value = '  -3.71%  '

# Row 37: 'Frax' -> 'MXToken'
$ws.Range("B37").Value = 'MXToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue "D37" '3.118'
$ws.Range("E37").Value = '  -1.91%  '

# Row 38: 'MXToken' -> 'RenderToken'
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue "D38" '2.346'
$ws.Range("E38").Value = '  -11.47%  '

# Row 39: 'RenderToken' -> 'VeChain'
$ws.Range("B39").Value = 'VeChain'
$ws.Range("C39").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
Set-TextValue "D39" '0.01831'
$ws.Range("E39").Value = '  -5.79%  '

# Row 40: 'VeChain' -> 'TheSandbox'
$ws.Range("B40").Value = 'TheSandbox'
$ws.Range("C40").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue "D40" '0.5041'
$ws.Range("E40").Value = '  -17.76%  '

# Row 41: 'TheSandbox' -> 'TrustWalletToken'
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue "D41" '0.9385'
$ws.Range("E41").Value = '  -12.19%  '

# Row 42: 'TrustWalletToken' -> 'Quant'
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue "D42" '115.93'
$ws.Range("E42").Value = '  +0.50%  '

# Row 43: 'Quant' -> 'FraxShare'
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
Set-TextValue "D43" '6.165'
$ws.Range("E43").Value = '  -4.74%  '

# Row 44: 'FraxShare' -> 'Aptos'
$ws.Range("B44").Value = 'Aptos'
$ws.Range("C44").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue "D44" '7.819'
$ws.Range("E44").Value = '  -10.39%  '

# Row 45
Set-TextValue "D45" '0.9995'
$ws.Range("E45").Value = '  -0.09%  '

# Row 46: 'Aptos' -> 'Algorand'
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue "D46" '0.1350'
$ws.Range("E46").Value = '  -9.51%  '

# Row 47: 'Algorand' -> 'Decentraland'
$ws.Range("B47").Value = 'Decentraland'
$ws.Range("C47").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
Set-TextValue "D47" '0.4388'
$ws.Range("E47").Value = '  -17.04%  '

# Row 48: 'Decentraland' -> 'Elrond'
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue "D48" '36.16'
$ws.Range("E48").Value = '  -3.09%  '

# Row 49
Set-TextValue "D49" '9.217'
$ws.Range("E49").Value = '  -7.32%  '

# Row 50: 'Elrond' -> 'Cronos'
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue "D50" '0.05938'
$ws.Range("E50").Value = '  -2.01%  '

# Row 51: 'Cronos' -> 'NEARProtocol'
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
Set-TextValue "D51" '1.477'
$ws.Range("E51").Value = '  -10.22%  '
